$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text first so numeric-looking values (e.g. "5.430",
# "0.9993") keep their exact literal formatting instead of being parsed as numbers
# and losing trailing zeros / significant digits (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.484.19'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '1.912.43'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('D8').Value = '0.2845'
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('D9').Value = '0.06707'
$ws.Range('E9').Value = '  -3.09%  '
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').Value = '102.44'
$ws.Range('E11').Value = '  -2.96%  '
$ws.Range('D12').Value = '0.07705'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.918.07'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('D15').Value = '0.6717'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '271.69'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '30.528.00'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').Value = '0.9993'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '0.000007469'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('D21').Value = '5.430'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('D22').Value = '0.4638'
$ws.Range('E22').Value = '  -10.05%  '
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '6.314'
$ws.Range('E24').Value = '  -3.74%  '
$ws.Range('D25').Value = '9.402'
$ws.Range('E25').Value = '  -3.44%  '
$ws.Range('D26').Value = '166.64'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '19.35'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').Value = '2.066'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('D29').Value = '1.383'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('D32').Value = '1.513'
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('D33').Value = '4.238'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').Value = '0.04736'
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('D35').Value = '0.7275'
$ws.Range('E35').Value = '  -4.14%  '
$ws.Range('D36').Value = '1.111'
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('D37').Value = '2.718'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').Value = '0.01924'
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('D39').Value = '2.611'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('D40').Value = '6.293'
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('D41').Value = '74.88'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('D42').Value = '1.968'
$ws.Range('E42').Value = '  -5.86%  '
$ws.Range('D43').Value = '0.8579'
$ws.Range('E43').Value = '  -5.08%  '
$ws.Range('D44').Value = '104.86'
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('D45').Value = '0.4267'
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('D46').Value = '0.9991'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = '7.424'
$ws.Range('E47').Value = '  -4.67%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1203'
$ws.Range('E48').Value = '  -3.69%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '920.03'
$ws.Range('E49').Value = '  -7.29%  '
$ws.Range('D50').Value = '34.82'
$ws.Range('E50').Value = '  -3.71%  '
$ws.Range('D51').Value = '8.813'
$ws.Range('E51').Value = '  -5.27%  '

# Reset styling on the Price column back to the default (no explicit style index),
# matching the source file where D2:D51 carry no "s" attribute.
$ws.Range("D2:D51").Style = "Normal"
